$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-7: RoundNumber, Dragon Hand, Tiger Hand, Result
$data = @(
    @(1, "3 of Dimonds", "A of Spades", "Dragon Won"),
    @(1, "3 of Dimonds", "A of Spades", "Dragon Won"),
    @(2, "6 of Spades",  "J of Heart",  "Tiger Won"),
    @(1, "3 of Dimonds", "A of Spades", "Dragon Won"),
    @(2, "6 of Spades",  "J of Heart",  "Tiger Won"),
    @(3, "9 of Heart",   "4 of Spades", "Dragon Won")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}
